$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2)
}

Replace-Text "2023-10-13 Friday" "2023-10-14 Saturday"

Replace-Text "29×23=" "48×21="
Replace-Text "71×31=" "56×25="
Replace-Text "81×70=" "40×70="
Replace-Text "91×50=" "95×28="
Replace-Text "62×21=" "33×53="

Replace-Text "52×62=" "89×30="
Replace-Text "93×16=" "45×28="
Replace-Text "20×32=" "81×44="
Replace-Text "79×66=" "33×19="
Replace-Text "26×87=" "85×69="

Replace-Text "81×42=" "76×94="
Replace-Text "17×80=" "70×58="
Replace-Text "68×39=" "26×48="
Replace-Text "78×90=" "23×41="
Replace-Text "58×21=" "80×66="

Replace-Text "24×13=" "84×92="
Replace-Text "26×65=" "76×22="
Replace-Text "13×20=" "35×14="
Replace-Text "39×75=" "79×62="
Replace-Text "48×75=" "62×56="

Replace-Text "69×23=" "36×91="
Replace-Text "15×74=" "83×97="
Replace-Text "21×84=" "70×72="
Replace-Text "25×44=" "79×74="
Replace-Text "68×76=" "85×46="
